$d = $word.ActiveDocument

# Locate the paragraph that starts the "Create an Quarto file; title it ..."
# list item so we can replace its runs with the expanded instructions.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create an Quarto file; title it*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$pRange = $target.Range
# Range of just the paragraph's text, excluding the trailing paragraph mark,
# so the replacement keeps the paragraph's own pPr (style/numbering) intact.
$textRange = $d.Range($pRange.Start, $pRange.End - 1)

$newRunsXml = '<w:r><w:t xml:space="preserve">Create an Quarto file (</w:t></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">File&gt;&gt;New File&gt;&gt;Quarto Document</w:t></w:r><w:r><w:t xml:space="preserve">). Title it</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">&#8220;Week 3 Homework.qmd&#8221;</w:t></w:r><w:r><w:t xml:space="preserve">, select</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">&#8220;HTML&#8221;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">output, and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">uncheck</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">&#8220;use visual markdown editor&#8221;</w:t></w:r><w:r><w:t xml:space="preserve">. (You can switch to this later.)</w:t></w:r>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$textRange.InsertXML($packageXml)

Write-Output "Paragraph now reads: $($target.Range.Text)"
